$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.471.83"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.571.59"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'288.24"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("D7").Value = "'0.3719"
$ws.Range("E7").Value = "  +0.84%  "
$ws.Range("D8").Value = "'48.32"
$ws.Range("E8").Value = "  -3.72%  "
$ws.Range("D9").Value = "'0.3319"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "'1.134"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "'20.74"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Value = "'5.935"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "'6.894"
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").Value = "1.572.27"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "'0.00001119"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "'87.84"
$ws.Range("E18").Value = "  -2.59%  "
$ws.Range("D19").Value = "'0.06744"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "'6.355"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").Value = "'16.55"
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").Value = "'12.07"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "22.465.78"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "'2.386"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").Value = "'2.576"
$ws.Range("E26").Value = "  -2.44%  "
$ws.Range("D27").Value = "'153.17"
$ws.Range("E27").Value = "  +2.63%  "
$ws.Range("D28").Value = "'19.69"
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("D29").Value = "'5.015"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").Value = "'124.39"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "1.745.58"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "'1.055"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").Value = "'6.141"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").Value = "'9.789"
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("D36").Value = "'0.08335"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").Value = "'0.2270"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("D39").Value = "'0.06408"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("D41").Value = "'5.356"
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.6304"
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'11.29"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "'13.85"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("D46").Value = "'0.6150"
$ws.Range("E46").Value = "  +5.79%  "
$ws.Range("D47").Value = "'3.774"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "'2.057"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").Value = "'125.49"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").Value = "'1.211"
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("D51").Value = "'0.07223"
$ws.Range("E51").Value = "  -1.54%  "
